$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email addresses for the existing students (rows 2-6)
$ws.Range("C2").Value = "pasne.d@husky.neu.edu"
$ws.Range("C3").Value = "sood.s@husky.neu.edu"
$ws.Range("C4").Value = "shail@ccs.neu.edu"
$ws.Range("C5").Value = "dave.v@husky.neu.edu"
$ws.Range("C6").Value = "snow.j@husky.neu.edu"

# Replace the placeholder student names (rows 7-10) with real names
$ws.Range("B7").Value = "Danny"
$ws.Range("B8").Value = "Erica"
$ws.Range("B9").Value = "Flurry"
$ws.Range("B10").Value = "Gara"

# ... and give them their real email addresses
$ws.Range("C7").Value = "danny.d@husky.neu.edu"
$ws.Range("C8").Value = "sniper.e@husky.neu.edu"
$ws.Range("C9").Value = "majin.f@husky.neu.edu"
$ws.Range("C10").Value = "hawking.g@husky.neu.edu"

# Add explicit mailto hyperlinks for the newly-set email addresses.
# Hyperlinks.Add() re-applies the Hyperlink cell style as a freshly minted
# style record instead of reusing the existing one, so re-apply the named
# "Hyperlink" style afterwards to put the cell back on the original shared
# style (matches the rest of column C).
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:shail@ccs.neu.edu") | Out-Null
$ws.Range("C4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:snow.j@husky.neu.edu") | Out-Null
$ws.Range("C6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:danny.d@husky.neu.edu") | Out-Null
$ws.Range("C7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:sniper.e@husky.neu.edu") | Out-Null
$ws.Range("C8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:majin.f@husky.neu.edu") | Out-Null
$ws.Range("C9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:hawking.g@husky.neu.edu") | Out-Null
$ws.Range("C10").Style = "Hyperlink"

# Move the active selection to C16
$ws.Range("C16").Select() | Out-Null
